$wb = $excel.ActiveWorkbook

# --- Fix data values on PROJECTDATA sheet ---
$wsProject = $wb.Worksheets.Item("PROJECTDATA")
$wsProject.Range("C3").Value = "On Goging"
$wsProject.Range("C4").Value = "Completed"

# --- Make PROJECTDATA the active sheet with C4 selected ---
$wsProject.Activate()
$wsProject.Range("C4").Select()

# --- USERCREDENTIAL keeps its own selection at H6 (no longer the active tab) ---
$wsUser = $wb.Worksheets.Item("USERCREDENTIAL")
$wsUser.Range("H6").Select()

# --- Re-activate PROJECTDATA so it's the final active/visible tab ---
$wsProject.Activate()
